# Update automatico via Actualizar 02-07-2021 13-17-37
# Shifts the "Fecha" (column D) timestamps: a new check cycle's timestamp
# is recorded for rows 2-15, and the two older cycles shift down one slot
# (rows 16-29 take what used to be rows 2-15's value, rows 30-43 take what
# used to be rows 16-29's value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value = 44234.55383578015
$ws.Range("D16:D29").Value = 44234.53264789352
$ws.Range("D30:D43").Value = 44234.51146017361
